$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking row (row 11): Right marks per question 3 -> 5
$ws.Range("B11").Value = 5

# Update total row (row 12): Total right marks 57 -> 95
$ws.Range("B12").Value = 95

# Update correct/total marks text (E12): 54/84 -> 95/140
$ws.Range("E12").Value = "95/140"
